# Fix Bugs in Check package
# The worksheet has an empty column Q between the data in column P and the
# "n_size" helper column in column R. Delete the empty column Q so that the
# data in column R shifts left to become column Q.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daten_Beispiel")

$ws.Columns.Item(17).Delete()

# Update the view: scroll so column M is the left-most visible column and
# select the (now) last column Q, matching the post-edit selection state.
$ws.Application.Goto($ws.Range("M1"))
$ws.Range("Q1:Q1048576").Select()
